# Update countries & provincias Spain
#
# Refreshes the COVID "Pais" snapshot:
#  - bumps the "Datos actualizados ..." timestamp cell (A1)
#  - updates case counters for a handful of countries whose totals changed
#  - because three countries overtook their neighbours in total cases,
#    the affected rows are re-sorted (descending by "Casos totales"):
#    the country name together with its own statistics moves as a unit
#    to its new rank, so both column A and columns B:H are rewritten for
#    every row touched by the shuffle.
#
# $wb / $excel already resolve to the open workbook; grab the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 10:22"

# --- Austria (row 20): straight data refresh, rank unchanged ---
$ws.Range("B20").Value = 13962
$ws.Range("C20").Value = 17
$ws.Range("E20").Value = 6625

# --- Rumania (row 33): straight data refresh, rank unchanged ---
$ws.Range("E33").Value = 5130
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 318

# --- Filipinas overtakes Malasia and Mexico -> rows 37-39 re-sorted ---
$ws.Range("A37").Value = "Filipinas"
$ws.Range("B37").Value = 4932
$ws.Range("C37").Value = 284
$ws.Range("D37").Value = 242
$ws.Range("E37").Value = 4375
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 18
$ws.Range("H37").Value = 315

$ws.Range("A38").Value = "Malasia"
$ws.Range("B38").Value = 4683
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 2108
$ws.Range("E38").Value = 2499
$ws.Range("F38").Value = 66
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 76

$ws.Range("A39").Value = "Mexico"
$ws.Range("B39").Value = 4661
$ws.Range("C39").Value = 442
$ws.Range("D39").Value = 1843
$ws.Range("E39").Value = 2522
$ws.Range("F39").Value = 185
$ws.Range("G39").Value = 23
$ws.Range("H39").Value = 296

# --- Martinica overtakes Guatemala -> rows 122-123 swap ---
$ws.Range("A122").Value = "Martinica"
$ws.Range("B122").Value = 157
$ws.Range("C122").Value = 2
$ws.Range("D122").Value = 50
$ws.Range("E122").Value = 101
$ws.Range("F122").Value = 19
$ws.Range("H122").Value = 6

$ws.Range("A123").Value = "Guatemala"
$ws.Range("D123").Value = 19
$ws.Range("E123").Value = 131
$ws.Range("F123").Value = 3
$ws.Range("H123").Value = 5

# --- Tanzania overtakes Macao..Guam -> rows 151-161 shift down by one ---
$ws.Range("A151").Value = "Tanzania"
$ws.Range("B151").Value = 46
$ws.Range("C151").Value = 14
$ws.Range("D151").Value = 7
$ws.Range("E151").Value = 36
$ws.Range("F151").Value = 0
$ws.Range("H151").Value = 3

$ws.Range("A152").Value = "Macao"
$ws.Range("D152").Value = 10
$ws.Range("E152").Value = 35
$ws.Range("F152").Value = 1
$ws.Range("H152").Value = 0

$ws.Range("A153").Value = "Guyana"
$ws.Range("B153").Value = 45
$ws.Range("D153").Value = 8
$ws.Range("E153").Value = 31
$ws.Range("F153").Value = 3
$ws.Range("H153").Value = 6

$ws.Range("A154").Value = "Zambia"
$ws.Range("B154").Value = 43
$ws.Range("D154").Value = 30
$ws.Range("E154").Value = 11
$ws.Range("F154").Value = 1
$ws.Range("H154").Value = 2

$ws.Range("A155").Value = "Birmania"
$ws.Range("B155").Value = 41
$ws.Range("D155").Value = 2
$ws.Range("E155").Value = 35
$ws.Range("H155").Value = 4

$ws.Range("A156").Value = "Puerto Rico"
$ws.Range("B156").Value = 39
$ws.Range("D156").Value = 1
$ws.Range("E156").Value = 36
$ws.Range("H156").Value = 2

$ws.Range("A157").Value = "Guinea-Bisau"
$ws.Range("B157").Value = 38
$ws.Range("D157").Value = 0
$ws.Range("E157").Value = 38
$ws.Range("H157").Value = 0

$ws.Range("A158").Value = "Benin"
$ws.Range("B158").Value = 35
$ws.Range("D158").Value = 5
$ws.Range("E158").Value = 29
$ws.Range("H158").Value = 1

$ws.Range("A159").Value = "Eritrea"
$ws.Range("B159").Value = 34
$ws.Range("E159").Value = 34
$ws.Range("H159").Value = 0

$ws.Range("A160").Value = "Haiti"
$ws.Range("B160").Value = 33
$ws.Range("E160").Value = 30
$ws.Range("H160").Value = 3

$ws.Range("A161").Value = "Guam"
$ws.Range("D161").Value = 0
$ws.Range("E161").Value = 31
$ws.Range("H161").Value = 1
